# Swap the Opening Qty / Rate / Sales Qty / Sales Value values
# (columns B, E, F, G) between each of the following adjacent row pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(112, 113),
    @(232, 233),
    @(243, 244),
    @(364, 365),
    @(366, 367),
    @(375, 376),
    @(380, 381),
    @(382, 383),
    @(385, 386),
    @(473, 474),
    @(572, 573)
)

$cols = @("B", "E", "F", "G")

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
